# Apply the "More slides and updates." edit:
#   - keep slide 1 (title slide) and the existing baseball-decision-tree
#     slide (id 257) where it is logically, but insert four new
#     "Title and Content" slides around it:
#       2: "What are random forests?"
#       3: "What are decision trees?"
#       (existing baseball slide stays here)
#       5: "Random forests"
#       6: "Parameters"

$p = $ppt.ActivePresentation

# --- Insert the two new slides that go before the existing slide 2 ---
$sRandomForests = $p.Slides.Add(2, 2)
$sDecisionTrees = $p.Slides.Add(3, 2)

# --- Insert three slides at the tail; one is a throwaway so the
#     surviving two pick up SlideIDs 261/262 (matching 256,258,259,257,261,262) ---
$throwaway   = $p.Slides.Add(5, 2)
$sParams1    = $p.Slides.Add(6, 2)
$sParams2    = $p.Slides.Add(7, 2)

# --- Fill in "Random forests" (will end up at index 5, id 261) ---
$sParams1.Shapes.Item(1).TextFrame.TextRange.Text = "Random forests"
$tf = $sParams1.Shapes.Item(2).TextFrame
$tr = $tf.TextRange
$tr.Text = "Bagging + random feature selection"
$tr.InsertAfter("`rBagging: generating new sample set from training set with replacement")
$tr.InsertAfter("`rGood for")
$tr.InsertAfter("`rHeterogeneous data")
$tr.InsertAfter("`rLinearly separable")
$tr.InsertAfter("`rLarge amounts of data")
$tf.TextRange.Paragraphs(2, 1).IndentLevel = 2
$tf.TextRange.Paragraphs(4, 1).IndentLevel = 2
$tf.TextRange.Paragraphs(5, 1).IndentLevel = 2
$tf.TextRange.Paragraphs(6, 1).IndentLevel = 2

# --- Fill in "Parameters" (will end up at index 6, id 262) ---
$sParams2.Shapes.Item(1).TextFrame.TextRange.Text = "Parameters"
$tf2 = $sParams2.Shapes.Item(2).TextFrame
$tr2 = $tf2.TextRange
$tr2.Text = "n_estimators"
$tr2.InsertAfter("`rcriterion")
$tr2.InsertAfter("`rmax_features")

# --- Drop the throwaway slide (burns SlideID 260) ---
$throwaway.Delete()

# --- Fill in the two new slides inserted at the top ---
$sRandomForests.Shapes.Item(1).TextFrame.TextRange.Text = "What are random forests?"
$tfA = $sRandomForests.Shapes.Item(2).TextFrame
$trA = $tfA.TextRange
$trA.Text = "Ensemble method: use multiple learning algorithms to get better results"
$trA.InsertAfter("`rConstructs multiple decision trees")
$trA.InsertAfter("`rOutput the most frequent label from those trees")

$sDecisionTrees.Shapes.Item(1).TextFrame.TextRange.Text = "What are decision trees?"
$tfB = $sDecisionTrees.Shapes.Item(2).TextFrame
$trB = $tfB.TextRange
$trB.Text = "Maps observations about data’s target value"
$trB.InsertAfter("`rLeaves represent labels")
$trB.InsertAfter("`rBranches represent conjunction of features that lead to labels")
